# Refresh crypto ranking snapshot (GitHub Actions scheduled update).
# Updates Price (D) / Volume(1h) (E) per coin row, and swaps rows whose
# ranking order changed (EthereumClassic/Stellar, Bittensor/Kaspa).
# Numeric-looking Price strings are entered with a leading apostrophe so
# Excel keeps them as text (matching the source data's inline-string type)
# instead of silently converting them to numbers and dropping formatting
# such as trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.570.93"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").Value = "3.144.58"
$ws.Range("E3").Value = "  +1.04%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'238.24"
$ws.Range("E5").Value = "  +8.33%  "

$ws.Range("D6").Value = "'645.78"
$ws.Range("E6").Value = "  +3.73%  "

$ws.Range("E7").Value = "  +11.41%  "

$ws.Range("E8").Value = "  -5.26%  "

$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").Value = "3.140.99"
$ws.Range("E10").Value = "  +1.04%  "

$ws.Range("D11").Value = "'0.721"
$ws.Range("E11").Value = "  +0.23%  "

$ws.Range("E12").Value = "  +4.31%  "

$ws.Range("D13").Value = "'36.64"
$ws.Range("E13").Value = "  +6.38%  "

$ws.Range("E14").Value = "  -3.78%  "

$ws.Range("D15").Value = "'5.64"
$ws.Range("E15").Value = "  +4.65%  "

$ws.Range("D16").Value = "90.182.42"
$ws.Range("E16").Value = "  -1.06%  "

$ws.Range("D17").Value = "3.709.41"
$ws.Range("E17").Value = "  +0.49%  "

$ws.Range("D18").Value = "3.133.07"
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("D19").Value = "'3.72"
$ws.Range("E19").Value = "  -0.36%  "

$ws.Range("D20").Value = "'14.50"
$ws.Range("E20").Value = "  +2.91%  "

$ws.Range("D21").Value = "'0.0000215"
$ws.Range("E21").Value = "  -2.46%  "

$ws.Range("D22").Value = "'451.10"
$ws.Range("E22").Value = "  +3.28%  "

$ws.Range("D23").Value = "'5.70"
$ws.Range("E23").Value = "  +10.22%  "

$ws.Range("D24").Value = "'9.07"
$ws.Range("E24").Value = "  +3.11%  "

$ws.Range("D25").Value = "'6.05"
$ws.Range("E25").Value = "  -1.59%  "

$ws.Range("D26").Value = "'91.72"
$ws.Range("E26").Value = "  +4.82%  "

$ws.Range("D27").Value = "'12.49"
$ws.Range("E27").Value = "  +2.40%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").Value = "'9.90"
$ws.Range("E30").Value = "  +8.72%  "

$ws.Range("E31").Value = "  -3.54%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'27.34"
$ws.Range("E32").Value = "  +15.40%  "

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "'0.201"
$ws.Range("E33").Value = "  +33.38%  "

$ws.Range("D34").Value = "'3.89"
$ws.Range("E34").Value = "  +4.05%  "

$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "'518.75"
$ws.Range("E35").Value = "  -1.50%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.151"
$ws.Range("E36").Value = "  +6.28%  "

$ws.Range("E37").Value = "  +5.73%  "

$ws.Range("D38").Value = "'7.15"
$ws.Range("E38").Value = "  +0.90%  "

$ws.Range("E39").Value = "  +3.08%  "

$ws.Range("D40").Value = "'0.425"
$ws.Range("E40").Value = "  +10.24%  "

$ws.Range("D41").Value = "'0.0867"
$ws.Range("E41").Value = "  -1.05%  "

$ws.Range("D42").Value = "'22.21"
$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("E44").Value = "  -16.34%  "

$ws.Range("D45").Value = "'3.33"
$ws.Range("E45").Value = "  +40.29%  "

$ws.Range("D46").Value = "'1.95"
$ws.Range("E46").Value = "  +1.82%  "

$ws.Range("D47").Value = "'0.705"
$ws.Range("E47").Value = "  +13.97%  "

$ws.Range("D48").Value = "'149.84"
$ws.Range("E48").Value = "  +2.01%  "

$ws.Range("D49").Value = "'4.61"
$ws.Range("E49").Value = "  +10.27%  "

$ws.Range("D50").Value = "'45.55"
$ws.Range("E50").Value = "  +3.51%  "

$ws.Range("D51").Value = "'1.36"
$ws.Range("E51").Value = "  +4.60%  "
